$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Correct the timestamp in A83: drop the intraday time component so it
#    matches the normalized 07:00 timestamp used by the other rows.
$ws.Range("A83").Value = 45462.2916666667

# 2. Append the new day's OHLCV/ticker data as row 84.
$ws.Range("A84").Value = 45463.4620486111
$ws.Range("B84").Value = 3000
$ws.Range("C84").Value = 2.95000004768372
$ws.Range("D84").Value = 2.91000008583069
$ws.Range("E84").Value = 2.91000008583069
$ws.Range("F84").Value = 2.95000004768372

# Re-apply the exact cell style already used elsewhere in the column so no
# new (duplicate) style entries are introduced in the workbook.
$ws.Range("A83").Copy()
$ws.Range("A84").PasteSpecial(-4122)
$ws.Range("B83").Copy()
$ws.Range("B84").PasteSpecial(-4122)
$ws.Range("C83").Copy()
$ws.Range("C84").PasteSpecial(-4122)
$ws.Range("D83").Copy()
$ws.Range("D84").PasteSpecial(-4122)
$ws.Range("E83").Copy()
$ws.Range("E84").PasteSpecial(-4122)
$ws.Range("F83").Copy()
$ws.Range("F84").PasteSpecial(-4122)

# G84/H84 hold text (the string representation of the close price, and the
# ticker) rather than numeric values. Copy whole cells that already contain
# the exact same text elsewhere in the sheet so the existing shared-string
# entries are reused and no cell is mis-typed as a number.
$ws.Range("G82").Copy()
$ws.Range("G84").PasteSpecial(-4163)
$ws.Range("H83").Copy()
$ws.Range("H84").PasteSpecial(-4163)
